$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web")
$ws.Activate()

# Row 228: new sub-section header "BookStore" (ElementID column only)
$ws.Range("A228").Value = "BookStore"

# Row 229: bookStoreSearch locator
$ws.Range("A229").Value = "bookStoreSearch"
$ws.Range("B229").Value = "//input[@id='searchBox']"
$ws.Range("C229").Value = "By.xpath"

# Row 230: selectBook locator
$ws.Range("A230").Value = "selectBook"
$ws.Range("B230").Value = "//*[@id=""see-book-Eloquent JavaScript, Second Edition""]/a"
$ws.Range("C230").Value = "By.xpath"

# Match column C's existing formatting used throughout the sheet
$ws.Range("C227").Copy() | Out-Null
$ws.Range("C229:C230").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("C229").Value = "By.xpath"
$ws.Range("C230").Value = "By.xpath"

# Update the view state to match where the user ended up after the edit
$ws.Range("A232").Select()
$window = $excel.ActiveWindow
$window.ScrollRow = 223
